$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (order matters for shared-string table indices)
$ws.Range("N1").Value = "Fitted Ratio"
$ws.Range("L1").Value = "Reverse Ratio"
$ws.Range("O1").Value = "Fitted Reverse Ratio"

# L column: Reverse Ratio = 1/K
$ws.Range("L2").Formula = "=1/K2"
$ws.Range("L3:L13").Formula = "=1/K3"
$ws.Range("L2:L13").NumberFormat = "0.00E+00"

# N column: Fitted Ratio (polynomial fit vs D)
$ws.Range("N2").Formula = "=(D2^2)*-0.014887522+0.420918*D2+0.9383"
$ws.Range("N3:N13").Formula = "=(D3^2)*-0.014887522+0.420918*D3+0.9383"

# O column: Fitted Reverse Ratio (polynomial fit vs D)
$ws.Range("O2").Formula = "=0.0112817207083844*(D2^2) +D2*( -0.408181546694434) + 3.7160769331443"
$ws.Range("O3:O13").Formula = "=0.0112817207083844*(D3^2) +D3*( -0.408181546694434) + 3.7160769331443"

# Update active selection
$null = $ws.Range("P8").Select()
